# Feria Lagunitas de Puerto Montt - Pepino ensalada
# A new weekly price observation (2022-10-07) is inserted as a new data
# row right after the existing 2022-02-22 entry (row 260), pushing every
# following row down by one. The sheet's used range therefore grows from
# A1:R331 to A1:R332.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 260, shifting rows 260:331 down to 261:332.
$ws.Range("A260").EntireRow.Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(260, 1).Value  = 4
$ws.Cells.Item(260, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(260, 3).Value  = "Los Lagos"
$ws.Cells.Item(260, 4).Value  = 44841
$ws.Cells.Item(260, 5).Value  = 10
$ws.Cells.Item(260, 6).Value  = 100112043
$ws.Cells.Item(260, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(260, 8).Value  = "Sin especificar"
$ws.Cells.Item(260, 9).Value  = "Primera"
$ws.Cells.Item(260, 10).Value = 400
$ws.Cells.Item(260, 11).Value = 32000
$ws.Cells.Item(260, 12).Value = 32000
$ws.Cells.Item(260, 13).Value = 32000
$ws.Cells.Item(260, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(260, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(260, 16).Value = 533
$ws.Cells.Item(260, 17).Value = 60
$ws.Cells.Item(260, 18).Value = "Hortaliza"
